# Updates cryptos list values (prices in column D, volume % change in column E).
# Rows 30/31 and 41/42 also swap coin name/link (B/C) per upstream reorder.
# D-column values that look like plain numbers ("215.52", "1.007", ...) are
# written with a leading apostrophe so Excel stores them as text (matching the
# original inline-string cells) instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.037.66'
$ws.Range("E2").Value = '  -2.17%  '

$ws.Range("D3").Value = '1.640.82'
$ws.Range("E3").Value = '  -2.20%  '

$ws.Range("E4").Value = '  -0.27%  '

$ws.Range("D5").Value = "'" + '215.52'
$ws.Range("E5").Value = '  -2.41%  '

$ws.Range("D6").Value = "'" + '0.5059'
$ws.Range("E6").Value = '  -2.94%  '

$ws.Range("E7").Value = '  -0.13%  '

$ws.Range("E8").Value = '  -0.32%  '

$ws.Range("D9").Value = "'" + '0.06406'
$ws.Range("E9").Value = '  -2.04%  '

$ws.Range("E10").Value = '  -3.38%  '

$ws.Range("E11").Value = '  +0.44%  '

$ws.Range("D12").Value = '1.652.91'
$ws.Range("E12").Value = '  -1.68%  '

$ws.Range("D13").Value = "'" + '4.248'
$ws.Range("E13").Value = '  -2.42%  '

$ws.Range("D14").Value = '1.867.44'
$ws.Range("E14").Value = '  -2.25%  '

$ws.Range("D15").Value = "'" + '0.5458'
$ws.Range("E15").Value = '  -3.00%  '

$ws.Range("D16").Value = '0.0₅7950'
$ws.Range("E16").Value = '  -1.35%  '

$ws.Range("D17").Value = "'" + '63.69'
$ws.Range("E17").Value = '  -2.91%  '

$ws.Range("D18").Value = '26.034.71'
$ws.Range("E18").Value = '  -2.40%  '

$ws.Range("D19").Value = "'" + '1.007'
$ws.Range("E19").Value = '  -0.08%  '

$ws.Range("D20").Value = "'" + '206.17'
$ws.Range("E20").Value = '  -3.89%  '

$ws.Range("D21").Value = "'" + '4.354'
$ws.Range("E21").Value = '  -3.59%  '

$ws.Range("D22").Value = "'" + '9.987'
$ws.Range("E22").Value = '  -1.91%  '

$ws.Range("D23").Value = "'" + '5.983'
$ws.Range("E23").Value = '  +0.78%  '

$ws.Range("E24").Value = '  -0.16%  '

$ws.Range("E25").Value = '  +11.73%  '

$ws.Range("D26").Value = "'" + '142.32'
$ws.Range("E26").Value = '  -0.96%  '

$ws.Range("E27").Value = '  -1.34%  '

$ws.Range("D28").Value = "'" + '6.874'
$ws.Range("E28").Value = '  -2.93%  '

$ws.Range("E29").Value = '  -0.75%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = "'" + '1.238'
$ws.Range("E30").Value = '  -2.87%  '

$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").Value = "'" + '0.05002'
$ws.Range("E31").Value = '  -4.98%  '

$ws.Range("D32").Value = "'" + '3.284'
$ws.Range("E32").Value = '  -2.97%  '

$ws.Range("D33").Value = "'" + '3.207'
$ws.Range("E33").Value = '  -1.63%  '

$ws.Range("D34").Value = "'" + '1.538'
$ws.Range("E34").Value = '  -4.20%  '

$ws.Range("D35").Value = "'" + '2.340'

$ws.Range("D36").Value = "'" + '0.9096'
$ws.Range("E36").Value = '  -2.54%  '

$ws.Range("D37").Value = "'" + '2.644'
$ws.Range("E37").Value = '  -5.04%  '

$ws.Range("D38").Value = "'" + '0.5681'
$ws.Range("E38").Value = '  -1.53%  '

$ws.Range("D39").Value = '1.125.95'
$ws.Range("E39").Value = '  -3.42%  '

$ws.Range("D40").Value = "'" + '0.01562'
$ws.Range("E40").Value = '  -3.41%  '

$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").Value = "'" + '1.007'
$ws.Range("E41").Value = '  -0.09%  '

$ws.Range("B42").Value = 'mCoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range("D42").Value = "'" + '2.557'
$ws.Range("E42").Value = '  -0.49%  '

$ws.Range("D43").Value = "'" + '5.636'
$ws.Range("E43").Value = '  -1.60%  '

$ws.Range("D44").Value = "'" + '0.8147'
$ws.Range("E44").Value = '  -2.18%  '

$ws.Range("D45").Value = "'" + '99.76'
$ws.Range("E45").Value = '  -0.14%  '

$ws.Range("D46").Value = '1.781.50'
$ws.Range("E46").Value = '  -2.21%  '

$ws.Range("E47").Value = '  +1.93%  '

$ws.Range("D48").Value = "'" + '0.4530'
$ws.Range("E48").Value = '  +0.63%  '

$ws.Range("D49").Value = "'" + '1.007'
$ws.Range("E49").Value = '  +0.01%  '

$ws.Range("E50").Value = '  -2.31%  '

$ws.Range("D51").Value = "'" + '7.706'
$ws.Range("E51").Value = '  -3.61%  '
